# Atualização de projeto, ajustes.
# Applies the data-entry corrections made to TestData.xlsx:
#  - "Cadastro" sheet: B2 nickname value changed ("Wilkerbn22" -> "Wilkerbn504")
#  - "Produtos" sheet: A8 "999" re-entered as text (quote-prefixed) instead of a number,
#    and the "Posição" labels in column C (rows 8, 12, 13, 14) updated.

$wb = $excel.ActiveWorkbook

$wsCadastro = $wb.Worksheets.Item("Cadastro")
$wsCadastro.Range("B2").Value = "Wilkerbn504"

$wsProdutos = $wb.Worksheets.Item("Produtos")
# Re-enter 999 as text (leading apostrophe = quote-prefixed literal, as Excel does
# when a user types '999 into a General formatted cell).
$wsProdutos.Range("A8").Value = "'999"
$wsProdutos.Range("C8").Value = "Q(7,0) V(7,1)"
$wsProdutos.Range("C12").Value = "C(11,0) V(11,1)"
$wsProdutos.Range("C13").Value = "C(12,0) V(12,1)"
$wsProdutos.Range("C14").Value = "C(13,0) V(13,1)"

# Keep the on-screen selection in sync with the saved file (cosmetic, matches the
# recorded workbook view state).
$wsProdutos.Range("B18").Select() | Out-Null

Write-Output "edit complete"
